$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 25006704
$ws.Range("I62").Value2 = 1700.4286
$ws.Range("J62").Value2 = 83351710
$ws.Range("K62").Value2 = 1700.4286
$ws.Range("L62").Value2 = 83351710
$ws.Range("M62").Value2 = -1076.4286
$ws.Range("N62").Value2 = -83352958

$ws.Range("H65").Value2 = 25006704
$ws.Range("I65").Value2 = 1700.4286
$ws.Range("J65").Value2 = 83351710
$ws.Range("K65").Value2 = 8502.143
$ws.Range("L65").Value2 = 416758550
$ws.Range("M65").Value2 = -5382.143
$ws.Range("N65").Value2 = -416764790

$ws.Range("H116").Value2 = 6927.1
$ws.Range("I116").Value2 = 7970
$ws.Range("J116").Value2 = 6231.8335
$ws.Range("K116").Value2 = 7970
$ws.Range("L116").Value2 = 6231.8335
$ws.Range("M116").Value2 = -4528
$ws.Range("N116").Value2 = -13115.8335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1808.4333
$ws.Range("I2").Value2 = 1210.909
$ws.Range("K2").Value2 = 1210.909
$ws.Range("M2").Value2 = -1097.909

$ws.Range("H32").Value2 = 9625.33
$ws.Range("I32").Value2 = 4649.7607
$ws.Range("J32").Value2 = 21806.896
$ws.Range("K32").Value2 = 4649.7607
$ws.Range("L32").Value2 = 21806.896
$ws.Range("M32").Value2 = -4362.7607
$ws.Range("N32").Value2 = -22380.896

$ws.Range("H45").Value2 = 3429.5652
$ws.Range("I45").Value2 = 2689.375
$ws.Range("J45").Value2 = 5121.4287
$ws.Range("K45").Value2 = 2689.375
$ws.Range("L45").Value2 = 5121.4287
$ws.Range("M45").Value2 = -2312.375
$ws.Range("N45").Value2 = -5875.4287

$ws.Range("H102").Value2 = 125000584
$ws.Range("I102").Value2 = 669.8570999999999
$ws.Range("J102").Value2 = 1000000000
$ws.Range("K102").Value2 = 669.8570999999999
$ws.Range("L102").Value2 = 1000000000
$ws.Range("M102").Value2 = 952.1429000000001
$ws.Range("N102").Value2 = -1000003244

$ws.Range("H110").Value2 = 3179.5483
$ws.Range("I110").Value2 = 3541.4075
$ws.Range("J110").Value2 = 737
$ws.Range("K110").Value2 = 3541.4075
$ws.Range("L110").Value2 = 737
$ws.Range("M110").Value2 = -1496.4075
$ws.Range("N110").Value2 = -4827

$ws.Range("H116").Value2 = 1808.4333
$ws.Range("I116").Value2 = 1210.909
$ws.Range("K116").Value2 = 1210.909
$ws.Range("M116").Value2 = 1083.091

$ws.Range("H132").Value2 = 2620.28
$ws.Range("I132").Value2 = 2133.35
$ws.Range("J132").Value2 = 4568
$ws.Range("K132").Value2 = 6400.049999999999
$ws.Range("L132").Value2 = 13704
$ws.Range("M132").Value2 = -3870.049999999999
$ws.Range("N132").Value2 = -18764

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1808.4333
$ws.Range("I3").Value2 = 1210.909
$ws.Range("K3").Value2 = 1210.909
$ws.Range("M3").Value2 = -1096.909

$ws.Range("H107").Value2 = 3384.5186
$ws.Range("I107").Value2 = 2806.9473
$ws.Range("J107").Value2 = 4756.25
$ws.Range("K107").Value2 = 2806.9473
$ws.Range("L107").Value2 = 4756.25
$ws.Range("M107").Value2 = -886.9472999999998
$ws.Range("N107").Value2 = -8596.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 4533.3335
$ws.Range("I16").Value2 = 4533.3335
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 4533.3335
$ws.Range("L16").Value2 = 0
$ws.Range("M16").Value2 = -4246.3335
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value2 = 2067.3494
$ws.Range("I31").Value2 = 1667.7727
$ws.Range("J31").Value2 = 2518.1538
$ws.Range("K31").Value2 = 1667.7727
$ws.Range("L31").Value2 = 2518.1538
$ws.Range("M31").Value2 = -1372.7727
$ws.Range("N31").Value2 = -3108.1538

$ws.Range("H34").Value2 = 2067.3494
$ws.Range("I34").Value2 = 1667.7727
$ws.Range("J34").Value2 = 2518.1538
$ws.Range("K34").Value2 = 1667.7727
$ws.Range("L34").Value2 = 2518.1538
$ws.Range("M34").Value2 = -1465.7727
$ws.Range("N34").Value2 = -2922.1538

$ws.Range("H99").Value2 = 2209.1538
$ws.Range("I99").Value2 = 1243.8889
$ws.Range("J99").Value2 = 4381
$ws.Range("K99").Value2 = 1243.8889
$ws.Range("L99").Value2 = 4381
$ws.Range("M99").Value2 = 254.1111000000001
$ws.Range("N99").Value2 = -7377

$ws.Range("H107").Value2 = 705.36365
$ws.Range("I107").Value2 = 849.7273
$ws.Range("J107").Value2 = 416.63635
$ws.Range("K107").Value2 = 849.7273
$ws.Range("L107").Value2 = 416.63635
$ws.Range("M107").Value2 = 1070.2727
$ws.Range("N107").Value2 = -4256.63635

$ws.Range("H113").Value2 = 4533.3335
$ws.Range("I113").Value2 = 4533.3335
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 4533.3335
$ws.Range("L113").Value2 = 0
$ws.Range("M113").Value2 = -2363.3335
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value2 = 1587.875
$ws.Range("I122").Value2 = 982.125
$ws.Range("J122").Value2 = 2799.375
$ws.Range("K122").Value2 = 2946.375
$ws.Range("L122").Value2 = 8398.125
$ws.Range("M122").Value2 = -496.375
$ws.Range("N122").Value2 = -13298.125

$ws.Range("H126").Value2 = 2209.1538
$ws.Range("I126").Value2 = 1243.8889
$ws.Range("J126").Value2 = 4381
$ws.Range("K126").Value2 = 3731.6667
$ws.Range("L126").Value2 = 13143
$ws.Range("M126").Value2 = -1261.6667
$ws.Range("N126").Value2 = -18083

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value2 = 868.64
$ws.Range("I11").Value2 = 166.15384
$ws.Range("J11").Value2 = 1629.6666
$ws.Range("K11").Value2 = 498.46152
$ws.Range("L11").Value2 = 4888.9998
$ws.Range("M11").Value2 = -358.46152
$ws.Range("N11").Value2 = -5168.9998

$ws.Range("H22").Value2 = 1160.0769
$ws.Range("J22").Value2 = 1215
$ws.Range("L22").Value2 = 3645
$ws.Range("N22").Value2 = -3983

$ws.Range("H26").Value2 = 294.73685
$ws.Range("I26").Value2 = 300
$ws.Range("J26").Value2 = 292.85715
$ws.Range("K26").Value2 = 900
$ws.Range("L26").Value2 = 878.5714499999999
$ws.Range("M26").Value2 = -612
$ws.Range("N26").Value2 = -1454.57145

$ws.Range("H27").Value2 = 1160.0769
$ws.Range("J27").Value2 = 1215
$ws.Range("L27").Value2 = 3645
$ws.Range("N27").Value2 = -3849

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 1919.238
$ws.Range("I102").Value2 = 989
$ws.Range("J102").Value2 = 4896
$ws.Range("K102").Value2 = 989
$ws.Range("L102").Value2 = 4896
$ws.Range("M102").Value2 = 633
$ws.Range("N102").Value2 = -8140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 2644.5454
$ws.Range("I61").Value2 = 1598.5714
$ws.Range("J61").Value2 = 4475
$ws.Range("K61").Value2 = 1598.5714
$ws.Range("L61").Value2 = 4475
$ws.Range("M61").Value2 = -1396.5714
$ws.Range("N61").Value2 = -4879

$ws.Range("H113").Value2 = 2644.5454
$ws.Range("I113").Value2 = 1598.5714
$ws.Range("J113").Value2 = 4475
$ws.Range("K113").Value2 = 1598.5714
$ws.Range("L113").Value2 = 4475
$ws.Range("M113").Value2 = 571.4286
$ws.Range("N113").Value2 = -8815

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 1620.1351
$ws.Range("I107").Value2 = 957
$ws.Range("K107").Value2 = 2871
$ws.Range("M107").Value2 = -951

$ws.Range("H113").Value2 = 39241.848
$ws.Range("I113").Value2 = 45769.953
$ws.Range("J113").Value2 = 3337.25
$ws.Range("K113").Value2 = 137309.859
$ws.Range("L113").Value2 = 10011.75
$ws.Range("M113").Value2 = -135139.859
$ws.Range("N113").Value2 = -14351.75
